# Fix bullets/lists indentation (first level was slightly indented to the
# right instead of right on the margin).
#
# For each of the 9 levels of the abstract numbering definition used by the
# document's bullet list, drop the now-redundant <w:tabs><w:tab w:val="num"
# w:pos="..."/></w:tabs> element from the level's <w:pPr> and bump the
# paragraph's <w:ind w:left="..."/> by 240 twips (so the text lines up flush
# with the margin instead of being offset by the old tab stop).

$d = $word.ActiveDocument

# (old tab w:pos, old ind w:left) pairs present in the original numbering.xml
$levels = @(
    @{tab = 0;    left = 480},
    @{tab = 720;  left = 1200},
    @{tab = 1440; left = 1920},
    @{tab = 2160; left = 2640},
    @{tab = 2880; left = 3360},
    @{tab = 3600; left = 4080},
    @{tab = 4320; left = 4800},
    @{tab = 5040; left = 5520},
    @{tab = 5760; left = 6240}
)

$xml = $d.WordOpenXML

foreach ($lvl in $levels) {
    $tabPos = $lvl.tab
    $oldLeft = $lvl.left
    $newLeft = $oldLeft + 240

    $oldFragment = '<w:tabs><w:tab w:val="num" w:pos="' + $tabPos + '" /></w:tabs><w:ind w:left="' + $oldLeft + '" w:hanging="480" />'
    $newFragment = '<w:ind w:left="' + $newLeft + '" w:hanging="480" />'

    $xml = $xml -replace [regex]::Escape($oldFragment), $newFragment
}

$d.WordOpenXML = $xml

Write-Output "Updated numbering indentation for $($levels.Count) list levels"
